$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the four "score" helper columns (D,E,F,G) for rows 3..83.
# Row 2 already carries these formulas; rows 3..83 were left blank and are
# now being completed with the same formula pattern (as shared formulas,
# same as the two existing blocks already used for column N).

$ws.Range("D3:D66").Formula = "=IF((N3)>=50,10,"""")"
$ws.Range("E3:E66").Formula = "=IF((N3)>=50,20,"""")"
$ws.Range("F3:F66").Formula = "=IF((N3)>=50,20,"""")"
$ws.Range("G3:G66").Formula = "=IF((N3)>=51,IF((N3-50)>50,50,IF((N3-50)<0,0,(N3-50))), """" )"

$ws.Range("D67:D83").Formula = "=IF((N67)>=50,10,"""")"
$ws.Range("E67:E83").Formula = "=IF((N67)>=50,20,"""")"
$ws.Range("F67:F83").Formula = "=IF((N67)>=50,20,"""")"
$ws.Range("G67:G83").Formula = "=IF((N67)>=51,IF((N67-50)>50,50,IF((N67-50)<0,0,(N67-50))), """" )"

# --- A handful of rows gained newly-reported input scores, which ripple
# into column N (the total) and therefore into D/E/F/G above.
$ws.Range("J4").Value = 17
$ws.Range("L30").Value = 25
$ws.Range("L36").Value = 40
$ws.Range("J43").Value = 16
$ws.Range("L46").Value = 40
$ws.Range("J53").Value = 11
$ws.Range("J83").Value = 11

# --- Selection moved from the old F11 cursor position to the freshly
# completed D column.
$ws.Range("D2:D83").Select() | Out-Null
